$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Database Interaction section: fill in evaluation scores (were 0, now 4) ---
$ws.Range("D26").Value = 4
$ws.Range("D27").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("D29").Value = 4

# Add the new comment explaining the direct DB integration / jupyter notebook
$ws.Range("E26").Value = "Ho anche creato un jupyter notebook che permette l'inserimento, a piacimento dell'utente, dei dataset"

# --- Row 43 grew taller (reformatted markdown cell) ---
$ws.Rows.Item(43).RowHeight = 48

# --- Update the view: scroll down and move the selection ---
$ws.Activate()
$ws.Range("A14").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E33").Select()

Write-Host "Edits applied"
